$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('B2').Value = 'Ameaças cibernéticas e lucros cessantes continuam entre principais riscos empresariais no mundo'
$ws.Range('C2').Value = 'https://www.revistaapolice.com.br/2023/03/ameacas-ciberneticas-e-lucros-cessantes-continuam-entre-principais-riscos-empresariais-no-mundo/'
$ws.Range('D2').Value = 'Segundo um relatório da AGCS, catástrofes naturais e a mudança climática caem no ranking à medida que as empresas priorizam as preocupações macroeconômicas urgentes'

$ws.Range('B3').Value = 'Crescimento de roubo e furto de veículos reforça importância do seguro automóvel'
$ws.Range('C3').Value = 'https://www.revistaapolice.com.br/2023/04/crescimento-de-roubo-e-furto-de-veiculos-reforca-importancia-do-seguro-automovel/'
$ws.Range('D3').Value = 'De acordo com dados do IBGE, o Brasil registrou cerca de 64 roubos e furtos de veículos por hora em 2022'

$ws.Range('B4').Value = 'Inter Risk Services apresenta Dalve Ortolani como CCO'
$ws.Range('C4').Value = 'https://www.revistaapolice.com.br/2023/04/inter-risk-services-apresenta-dalve-ortolani-como-cco/'
$ws.Range('D4').Value = 'Além da contratação do executivo, a corretora reforçou a equipe nas operações do Rio de Janeiro, São Paulo e Salvador'

$ws.Range('B5').Value = 'Marcelo Blay e Manes Erlichman assumem nova posição na Creditas'
$ws.Range('C5').Value = 'https://www.revistaapolice.com.br/2023/04/marcelo-blay-e-manes-erlichman-assumem-nova-posicao-na-creditas/'
$ws.Range('D5').Value = 'Marcelo Blay e Manes Erlichman assumem a posição de senior advisors na Creditas, que adquiriu a Minuto Seguros em julho de 2021'

$ws.Range('B6').Value = 'Alessandro Octaviani toma posse como superintendente da Susep'
$ws.Range('C6').Value = 'https://www.revistaapolice.com.br/2023/04/alessandro-octaviani-toma-posse-como-superintendente-da-susep/'
$ws.Range('D6').Value = 'O novo superintendente entrará em exercício na próxima segunda-feira, 10 de abril, dedicando a primeira semana a reuniões e despachos internos, para apresentação dos projetos em curso na autarquia'

$ws.Range('B7').Value = 'Mercado de seguros aproveita crescimento da indústria pet para ofertar produtos'
$ws.Range('C7').Value = 'https://www.revistaapolice.com.br/2023/04/mercado-de-seguros-aproveita-crescimento-da-industria-pet-para-ofertar-produtos/'
$ws.Range('D7').Value = 'Levantamento do Instituto Pet Brasil prevê que o segmento fechou 2022 com um faturamento de R$ 59,9 bilhões'

$ws.Range('B8').Value = 'Encontro do CCS-SP discute retorno do corretor ao Open Insurance'
$ws.Range('C8').Value = 'https://www.revistaapolice.com.br/2023/04/encontro-do-ccs-sp-discute-retorno-do-corretor-de-seguros-ao-open-insurance/'
$ws.Range('D8').Value = 'Ex-diretor e ex-superintendente da Susep explicam aos associados detalhes da resolução que extinguiu a SISS e criou a SPOC'

$ws.Range('B9').Value = 'GBOEX aposta em novo conceito de comunicação'
$ws.Range('C9').Value = 'https://www.revistaapolice.com.br/2023/04/gboex-aposta-em-novo-conceito-de-comunicacao/'
$ws.Range('D9').Value = 'Empresa revisitou linha de comunicação e procurou atualizar projetos para traduzirem as alterações e inovações em sua gestão'

$ws.Range('B10').Value = 'Gestão financeira deve ser prioridade para pequenos empresários'
$ws.Range('C10').Value = 'https://www.revistaapolice.com.br/2023/04/gestao-financeira-deve-ser-prioridade-para-pequenos-empresarios/'
$ws.Range('D10').Value = 'Empreendedores devem apostar em organizar sua forma de controlar os gastos, lucros, investimentos e priorizar o uso da conta de banco jurídica'

$ws.Range('B11').Value = 'SulAmérica reforça linha Odonto PME e Empresarial'
$ws.Range('C11').Value = 'https://www.revistaapolice.com.br/2023/04/sulamerica-reforca-linha-odonto-pme-e-empresarial/'
$ws.Range('D11').Value = 'Após anunciar novidades no portfólio Odonto Individual, companhia apresenta lançamentos e novos benefícios para os planos PME, PME Mais e Empresarial'
